# Generate Report for Handoff
# Updates the localization-status report:
#  - Refreshes the "Latest HO Xliff Generate Date" (Overview) / "Latest
#    Handoff Datetime" (de-de) timestamps for the 544bb61c...fd84d22d
#    group of files from 00:19:37 -> 00:19:53.
#  - Refreshes the "Latest Handoff Datetime" timestamps on the zh-cn
#    sheet for the same group of files from 00:19:32 -> 00:19:48.
#  - Sets the "Priority" column to "ht" for that same group of files on
#    both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

for ($r = 8; $r -le 13; $r++) {
    # Overview: "Latest HO Xliff Generate Date" column (G)
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-25 00:19:53"

    # de-de: "Latest Handoff Datetime" column (H) shares the same text
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-25 00:19:53"

    # zh-cn: "Latest Handoff Datetime" column (H)
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-25 00:19:48"

    # Priority column (E) becomes "ht" on both zh-cn and de-de
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
}
